$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert 2 new columns before D (shifts old D..K to F..M) ---
$ws.Range("D:E").EntireColumn.Insert()

# --- Step 2: copy number formats from the shifted columns (F:G) into the
#     new D:E columns, limited to the row-blocks that actually hold data ---
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: populate the new D/E columns with the newest two quarters' data ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 798000
$ws.Range("E8").Value = 693000
$ws.Range("D9").Value = 279200
$ws.Range("E9").Value = 285400
$ws.Range("D10").Value = 518800
$ws.Range("E10").Value = 407600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 27400
$ws.Range("E14").Value = 3200
$ws.Range("D15").Value = 84300
$ws.Range("E15").Value = 79900
$ws.Range("D17").Value = 532700
$ws.Range("E17").Value = 503300
$ws.Range("D18").Value = 265300
$ws.Range("E18").Value = 189700
$ws.Range("D20").Value = -52100
$ws.Range("E20").Value = -54600
$ws.Range("D21").Value = 297400
$ws.Range("E21").Value = 214900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 213200
$ws.Range("E23").Value = 135100
$ws.Range("D24").Value = 58700
$ws.Range("E24").Value = 35200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 154500
$ws.Range("E26").Value = 99800
$ws.Range("D27").Value = 153100
$ws.Range("E27").Value = 100500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 52100
$ws.Range("E32").Value = 54600
$ws.Range("D33").Value = 153100
$ws.Range("E33").Value = 100500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 153100
$ws.Range("E35").Value = 100500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 145100
$ws.Range("E41").Value = 118400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 547300
$ws.Range("E43").Value = 551700
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 74700
$ws.Range("E45").Value = 96300
$ws.Range("D46").Value = 767100
$ws.Range("E46").Value = 766500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 731500
$ws.Range("E48").Value = 723900
$ws.Range("D49").Value = 5438100
$ws.Range("E49").Value = 5487100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 125300
$ws.Range("E52").Value = 141700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 7062000
$ws.Range("E54").Value = 7119100
$ws.Range("D57").Value = 67800
$ws.Range("E57").Value = 74200
$ws.Range("D58").Value = 96100
$ws.Range("E58").Value = 42100
$ws.Range("D59").Value = 240300
$ws.Range("E59").Value = 251000
$ws.Range("D60").Value = 404200
$ws.Range("E60").Value = 367300
$ws.Range("D61").Value = 3884900
$ws.Range("E61").Value = 4105900
$ws.Range("D62").Value = 904000
$ws.Range("E62").Value = 893900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5209300
$ws.Range("E66").Value = 5391000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 620400
$ws.Range("E72").Value = 484400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1852800
$ws.Range("E76").Value = 1728100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 153100
$ws.Range("E81").Value = 100500
$ws.Range("D83").Value = 70000
$ws.Range("E83").Value = 64800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 237200
$ws.Range("E89").Value = 175500
$ws.Range("D91").Value = -42700
$ws.Range("E91").Value = -27200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -13800
$ws.Range("E94").Value = -43700
$ws.Range("D96").Value = -17100
$ws.Range("E96").Value = -17100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -196700
$ws.Range("E100").Value = -161000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 26700
$ws.Range("E102").Value = -29300

# --- Step 4: a handful of prior-quarter cells were restated, not just shifted ---
$ws.Range("H17").Value = 517500
$ws.Range("I17").Value = 487200
$ws.Range("H18").Value = 136200
$ws.Range("I18").Value = 124700
$ws.Range("H20").Value = -48600
$ws.Range("I20").Value = -50600
$ws.Range("H32").Value = 48600
$ws.Range("I32").Value = 50600
$ws.Range("F91").Value = -15300
$ws.Range("G91").Value = -21100
$ws.Range("H91").Value = -23600
$ws.Range("I91").Value = -21200
$ws.Range("J91").Value = -13800
